$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 ("I0") and J1 ("IF"), matching H1's formatting ---
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2..37: column I (I0) and column J (IF) ---
# I0 values per row (defaults to 1, except rows 3 and 4)
$i0 = @{
    2 = 1;  3 = 7;  4 = 9;  5 = 1;  6 = 1;  7 = 1;  8 = 1;  9 = 1;  10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1; 20 = 1;
    21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1; 29 = 1; 30 = 1;
    31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 1; 36 = 1; 37 = 1
}

# IF values per row (equal to the existing IP/H value, except rows 3 and 4)
$if_ = @{
    2 = 5;  3 = 8;  4 = 9;  5 = 5;  6 = 4;  7 = 3;  8 = 6;  9 = 6;  10 = 6;
    11 = 7; 12 = 4; 13 = 4; 14 = 6; 15 = 7; 16 = 7; 17 = 6; 18 = 6; 19 = 5; 20 = 5;
    21 = 7; 22 = 8; 23 = 6; 24 = 5; 25 = 7; 26 = 5; 27 = 5; 28 = 3; 29 = 7; 30 = 6;
    31 = 6; 32 = 6; 33 = 5; 34 = 5; 35 = 3; 36 = 4; 37 = 2
}

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 9).Value = $i0[$r]
    $ws.Cells.Item($r, 10).Value = $if_[$r]
}
